$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A6").Value = 10
$ws.Range("B6").Value = 10

$ws.Range("D15").Select()
